$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$val) {
    # Force the cell to be stored as plain text (shared-string) even when the
    # text looks numeric (e.g. long phone numbers or "yyyy-mm-dd hh:mm:ss"
    # date strings), then reset the style back to "Normal" so the cell keeps
    # no style override, matching the formatting of the rest of the sheet.
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

# New records appended to the bottom of the table (rows 95-99).
$dates = @(
    "2026-02-04 11:13:24",
    "2026-02-20 15:13:09",
    "2026-02-20 06:44:25",
    "2026-02-09 00:40:08",
    "2026-02-16 08:11:49"
)
$numbers = @(
    "237675678961",
    "237681678622",
    "237654041671",
    "237675629624",
    "237683232376"
)
$names = @(
    "NOUPA KAMGAING AGNES CHIC MOBILE",
    "Marie Rosine Magne Talla",
    "LA NEGRESSE SARL FOKAM KOM DANICE KEVIN",
    "LA NEGRESSE LTDLA CBOX R0 MALLA TALLA JACQUELINE",
    "VAKENA SYLVIE YOK PASL (Pan African Saving and Loan)"
)
$balances = @(141, 49054, 49350, 14, 500214)

$startRow = 95
$count = $dates.Length

# Write column by column (all dates, then all numbers, then all names, then
# all balances) so new unique text values are registered in the same batch
# order they were introduced: Date column, Number column, Name column.
for ($i = 0; $i -lt $count; $i++) {
    Set-TextValue $ws.Cells.Item($startRow + $i, 1) $dates[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    Set-TextValue $ws.Cells.Item($startRow + $i, 2) $numbers[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    Set-TextValue $ws.Cells.Item($startRow + $i, 3) $names[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value2 = $balances[$i]
}
